$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.825.70'
$ws.Range('E2').Value = '  -0.93%  '
$ws.Range('D3').Value = '1.732.35'
$ws.Range('E3').Value = '  -1.83%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.0000'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '228.60'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -3.82%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.9998'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.08%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5230'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.62%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2747'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.09%  '
$ws.Range('E9').Value = '  -3.11%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.06120'
$ws.Range('D10').Style = 'Normal'
$ws.Range('D11').Value = '1.734.45'
$ws.Range('E11').Value = '  -1.73%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.07063'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.32%  '
$ws.Range('E13').Value = '  -6.92%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.6346'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -3.13%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '4.520'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.02%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '76.51'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -2.60%  '
$ws.Range('E17').Value = '  +0.19%  '
$ws.Range('E18').Value = '  +0.12%  '
$ws.Range('D19').Value = '25.817.22'
$ws.Range('E19').Value = '  -0.99%  '
$ws.Range('E20').Value = '  -2.44%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.000006626'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -1.87%  '
$ws.Range('D22').Value = '1.959.99'
$ws.Range('E22').Value = '  -1.81%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.195'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +2.29%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '8.767'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +3.82%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '5.161'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.95%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '140.04'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +1.57%  '
$ws.Range('E27').Value = '  +1.32%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '14.99'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -1.54%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.774'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -3.60%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '102.00'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -1.06%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.08270'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.89%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.706'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.29%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.498'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.82%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.04446'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.06%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.611'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -1.51%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.9675'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -3.83%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.6154'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.37%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.667'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -3.36%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01562'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -1.53%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.9998'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.07%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.899'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -4.36%  '
$ws.Range('B42').Value = 'Quant'
$ws.Range('C42').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '99.49'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -3.54%  '
$ws.Range('B43').Value = 'TheSandbox'
$ws.Range('C43').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.3808'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -2.78%  '
$ws.Range('B44').Value = 'FraxShare'
$ws.Range('C44').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '5.012'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.64%  '
$ws.Range('B45').Value = 'TrustWalletToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.7196'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -4.89%  '
$ws.Range('B46').Value = 'Cronos'
$ws.Range('C46').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.05330'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -3.19%  '
$ws.Range('B47').Value = 'Algorand'
$ws.Range('C47').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.1117'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -1.07%  '
$ws.Range('B48').Value = 'Aptos'
$ws.Range('C48').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '6.138'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -3.79%  '
$ws.Range('B49').Value = 'Aave'
$ws.Range('C49').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '53.14'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.12%  '
$ws.Range('B50').Value = 'Elrond'
$ws.Range('C50').Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '29.88'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -1.22%  '
$ws.Range('B51').Value = 'EnergySwap'
$ws.Range('C51').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '7.594'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +1.51%  '
